# Daily attendance processing - 2026-01-31 08:44:42
# Rotate the comma-separated "Recorded By" list in column G (right rotate
# by one place: the last entry moves to the front) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ($val -like "*, *") {
        $parts = $val -split ", "
        $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = $rotated -join ", "
    }
}
